# Sync attendance_reports: reorder the "Recorded By" (column G) entries so
# that multi-author lists are listed in reverse order, matching the
# canonical order synced from the main repo. A single exact combination
# ("dnasr281@gmail.com, admin@admin.com") is left untouched, exactly as in
# the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$exception = "dnasr281@gmail.com, admin@admin.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null) {
        if ($val -ne $exception) {
            if ($val.Contains(",")) {
                $rawParts = $val.Split(",")
                $parts = @()
                foreach ($p in $rawParts) {
                    $parts += $p.Trim()
                }
                [array]::Reverse($parts)
                $newVal = [string]::Join(", ", $parts)
                $cell.Value = $newVal
            }
        }
    }
}
